$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared-string based effic_choices values
$ws.Range("H3").Value = "[('Low 74%', 74), ('Standard 80%', 80), ('High Efficiency Condensing 95%', 95)]"
$ws.Range("H5").Value = "[('Low 74%', 74), ('Standard 80%', 80), ('High Efficiency (e.g. Toyostove) 84%', 84)]"
$ws.Range("H6").Value = "[('Low 74%', 74), ('Standard 80%', 80), ('High Efficiency 84%', 84)]"

# Update comment on G1
$comment = $ws.Range("G1").Comment
[void]$comment.Text("This is a typical DHW efficiency for the fuel type and is only used to back out DHW consumption from the user's entry of total fuel use if they state that the total includes DHW.")

# Update selection / view (scroll the view so column B is the leftmost visible
# column, then select H7 to match the saved cursor position)
[void]$ws.Range("H7").Select()
$excel.ActiveWindow.ScrollColumn = 2
